$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "BTC"
$ws.Range("C2").Value = "Bitcoin"
$ws.Range("D2").Value = 25800
$ws.Range("E2").Value = 499416847911
$ws.Range("F2").Value = 13388513728
$ws.Range("G2").Value = 0.41741
$ws.Range("B3").Value = "ETH"
$ws.Range("C3").Value = "Ethereum"
$ws.Range("D3").Value = 1753.6
$ws.Range("E3").Value = 210295110237
$ws.Range("F3").Value = 6727715091
$ws.Range("G3").Value = 0.57262
$ws.Range("B4").Value = "USDT"
$ws.Range("C4").Value = "Tether"
$ws.Range("D4").Value = 0.999409
$ws.Range("E4").Value = 83375403179
$ws.Range("F4").Value = 14960408295
$ws.Range("G4").Value = -0.16715
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "BNB"
$ws.Range("D5").Value = 236.92
$ws.Range("E5").Value = 36749154092
$ws.Range("F5").Value = 770548238
$ws.Range("G5").Value = -0.33121
$ws.Range("B6").Value = "USDC"
$ws.Range("C6").Value = "USD Coin"
$ws.Range("D6").Value = 1.001
$ws.Range("E6").Value = 28372340295
$ws.Range("F6").Value = 2897778863
$ws.Range("G6").Value = 0.05255
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "XRP"
$ws.Range("D7").Value = 0.509057
$ws.Range("E7").Value = 26416565691
$ws.Range("F7").Value = 1085072231
$ws.Range("G7").Value = 2.51725
$ws.Range("B8").Value = "STETH"
$ws.Range("C8").Value = "Lido Staked Ether"
$ws.Range("D8").Value = 1750.66
$ws.Range("E8").Value = 12529452955
$ws.Range("F8").Value = 23360654
$ws.Range("G8").Value = 0.45745
$ws.Range("B9").Value = "ADA"
$ws.Range("C9").Value = "Cardano"
$ws.Range("D9").Value = 0.272429
$ws.Range("E9").Value = 9407401686
$ws.Range("F9").Value = 786656730
$ws.Range("G9").Value = 13.35538
$ws.Range("B10").Value = "DOGE"
$ws.Range("C10").Value = "Dogecoin"
$ws.Range("D10").Value = 0.062361
$ws.Range("E10").Value = 8680589754
$ws.Range("F10").Value = 525062297
$ws.Range("G10").Value = 4.27137
$ws.Range("B11").Value = "TRX"
$ws.Range("C11").Value = "TRON"
$ws.Range("D11").Value = 0.069425
$ws.Range("E11").Value = 6239840447
$ws.Range("F11").Value = 316607938
$ws.Range("G11").Value = 1.30826
$ws.Range("B12").Value = "SOL"
$ws.Range("C12").Value = "Solana"
$ws.Range("D12").Value = 15.62
$ws.Range("E12").Value = 6159449074
$ws.Range("F12").Value = 750045136
$ws.Range("G12").Value = 6.46173
$ws.Range("B13").Value = "LTC"
$ws.Range("C13").Value = "Litecoin"
$ws.Range("D13").Value = 78.57
$ws.Range("E13").Value = 5720934564
$ws.Range("F13").Value = 678256531
$ws.Range("G13").Value = 1.82913
$ws.Range("B14").Value = "MATIC"
$ws.Range("C14").Value = "Polygon"
$ws.Range("D14").Value = 0.620169
$ws.Range("E14").Value = 5675193407
$ws.Range("F14").Value = 709691521
$ws.Range("G14").Value = 7.01136
$ws.Range("B15").Value = "DOT"
$ws.Range("C15").Value = "Polkadot"
$ws.Range("D15").Value = 4.49
$ws.Range("E15").Value = 5554238961
$ws.Range("F15").Value = 161987008
$ws.Range("G15").Value = 0.76047
$ws.Range("B16").Value = "BUSD"
$ws.Range("C16").Value = "Binance USD"
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 4780982453
$ws.Range("F16").Value = 1037599624
$ws.Range("G16").Value = -0.00883
$ws.Range("B17").Value = "DAI"
$ws.Range("C17").Value = "Dai"
$ws.Range("D17").Value = 0.999591
$ws.Range("E17").Value = 4541831427
$ws.Range("F17").Value = 83068865
$ws.Range("G17").Value = -0.00006
$ws.Range("B18").Value = "WBTC"
$ws.Range("C18").Value = "Wrapped Bitcoin"
$ws.Range("D18").Value = 25780
$ws.Range("E18").Value = 4031783326
$ws.Range("F18").Value = 89789242
$ws.Range("G18").Value = 0.29723
$ws.Range("B19").Value = "AVAX"
$ws.Range("C19").Value = "Avalanche"
$ws.Range("D19").Value = 11.69
$ws.Range("E19").Value = 4016941434
$ws.Range("F19").Value = 230145131
$ws.Range("G19").Value = 2.02142
$ws.Range("B20").Value = "SHIB"
$ws.Range("C20").Value = "Shiba Inu"
$ws.Range("D20").Value = 0.00000675
$ws.Range("E20").Value = 3954452367
$ws.Range("F20").Value = 183256117
$ws.Range("G20").Value = 4.57203
$ws.Range("B21").Value = "LEO"
$ws.Range("C21").Value = "LEO Token"
$ws.Range("D21").Value = 3.55
$ws.Range("E21").Value = 3284236364
$ws.Range("F21").Value = 544951
$ws.Range("G21").Value = 0.35146
$ws.Range("B22").Value = "UNI"
$ws.Range("C22").Value = "Uniswap"
$ws.Range("D22").Value = 4.06
$ws.Range("E22").Value = 3042886767
$ws.Range("F22").Value = 61679498
$ws.Range("G22").Value = 2.51607
$ws.Range("B23").Value = "LINK"
$ws.Range("C23").Value = "Chainlink"
$ws.Range("D23").Value = 5.19
$ws.Range("E23").Value = 2669744581
$ws.Range("F23").Value = 271578166
$ws.Range("G23").Value = 3.40861
$ws.Range("B24").Value = "XMR"
$ws.Range("C24").Value = "Monero"
$ws.Range("D24").Value = 137.22
$ws.Range("E24").Value = 2488575691
$ws.Range("F24").Value = 59530785
$ws.Range("G24").Value = 0.79191
$ws.Range("B25").Value = "OKB"
$ws.Range("C25").Value = "OKB"
$ws.Range("D25").Value = 40.86
$ws.Range("E25").Value = 2441105339
$ws.Range("F25").Value = 9029561
$ws.Range("G25").Value = -1.81827
$ws.Range("B26").Value = "ATOM"
$ws.Range("C26").Value = "Cosmos Hub"
$ws.Range("D26").Value = 8.26
$ws.Range("E26").Value = 2403508443
$ws.Range("F26").Value = 104690418
$ws.Range("G26").Value = 5.95728
$ws.Range("B27").Value = "XLM"
$ws.Range("C27").Value = "Stellar"
$ws.Range("D27").Value = 0.082961
$ws.Range("E27").Value = 2227062683
$ws.Range("F27").Value = 59405381
$ws.Range("G27").Value = 2.003
$ws.Range("B28").Value = "TON"
$ws.Range("C28").Value = "Toncoin"
$ws.Range("D28").Value = 1.46
$ws.Range("E28").Value = 2137523261
$ws.Range("F28").Value = 8168681
$ws.Range("G28").Value = -0.98084
$ws.Range("B29").Value = "ETC"
$ws.Range("C29").Value = "Ethereum Classic"
$ws.Range("D29").Value = 15.17
$ws.Range("E29").Value = 2137317550
$ws.Range("F29").Value = 104900932
$ws.Range("G29").Value = 4.60761
$ws.Range("B30").Value = "TUSD"
$ws.Range("C30").Value = "TrueUSD"
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 2038543739
$ws.Range("F30").Value = 1628220782
$ws.Range("G30").Value = 0.09704
$ws.Range("B31").Value = "BCH"
$ws.Range("C31").Value = "Bitcoin Cash"
$ws.Range("D31").Value = 102.81
$ws.Range("E31").Value = 1981960682
$ws.Range("F31").Value = 75512503
$ws.Range("G31").Value = 2.11127
$ws.Range("B32").Value = "ICP"
$ws.Range("C32").Value = "Internet Computer"
$ws.Range("D32").Value = 3.75
$ws.Range("E32").Value = 1622521313
$ws.Range("F32").Value = 32291673
$ws.Range("G32").Value = -1.03993
$ws.Range("B33").Value = "LDO"
$ws.Range("C33").Value = "Lido DAO"
$ws.Range("D33").Value = 1.8
$ws.Range("E33").Value = 1567446717
$ws.Range("F33").Value = 49908591
$ws.Range("G33").Value = -1.8169
$ws.Range("B34").Value = "QNT"
$ws.Range("C34").Value = "Quant"
$ws.Range("D34").Value = 102.07
$ws.Range("E34").Value = 1483160800
$ws.Range("F34").Value = 21260584
$ws.Range("G34").Value = -0.83519
$ws.Range("B35").Value = "FIL"
$ws.Range("C35").Value = "Filecoin"
$ws.Range("D35").Value = 3.44
$ws.Range("E35").Value = 1464954134
$ws.Range("F35").Value = 177236276
$ws.Range("G35").Value = 2.38016
$ws.Range("B36").Value = "CRO"
$ws.Range("C36").Value = "Cronos"
$ws.Range("D36").Value = 0.055206
$ws.Range("E36").Value = 1438248305
$ws.Range("F36").Value = 10516906
$ws.Range("G36").Value = 7.40031
$ws.Range("B37").Value = "HBAR"
$ws.Range("C37").Value = "Hedera"
$ws.Range("D37").Value = 0.04388459
$ws.Range("E37").Value = 1381268780
$ws.Range("F37").Value = 33749205
$ws.Range("G37").Value = 0.28063
$ws.Range("B38").Value = "ARB"
$ws.Range("C38").Value = "Arbitrum"
$ws.Range("D38").Value = 1.008
$ws.Range("E38").Value = 1276930728
$ws.Range("F38").Value = 311740124
$ws.Range("G38").Value = -1.21562
$ws.Range("B39").Value = "APT"
$ws.Range("C39").Value = "Aptos"
$ws.Range("D39").Value = 5.97
$ws.Range("E39").Value = 1189988763
$ws.Range("F39").Value = 93013984
$ws.Range("G39").Value = -0.05791
$ws.Range("B40").Value = "VET"
$ws.Range("C40").Value = "VeChain"
$ws.Range("D40").Value = 0.0155913
$ws.Range("E40").Value = 1127965332
$ws.Range("F40").Value = 53737517
$ws.Range("G40").Value = 4.80832
$ws.Range("B41").Value = "NEAR"
$ws.Range("C41").Value = "NEAR Protocol"
$ws.Range("D41").Value = 1.21
$ws.Range("E41").Value = 1098208500
$ws.Range("F41").Value = 75773337
$ws.Range("G41").Value = -0.62973
$ws.Range("B42").Value = "USDP"
$ws.Range("C42").Value = "Pax Dollar"
$ws.Range("D42").Value = 0.999436
$ws.Range("E42").Value = 1004307950
$ws.Range("F42").Value = 2064877
$ws.Range("G42").Value = -0.11238
$ws.Range("B43").Value = "FRAX"
$ws.Range("C43").Value = "Frax"
$ws.Range("D43").Value = 0.999644
$ws.Range("E43").Value = 1002394014
$ws.Range("F43").Value = 8771050
$ws.Range("G43").Value = -0.08224
$ws.Range("B44").Value = "GRT"
$ws.Range("C44").Value = "The Graph"
$ws.Range("D44").Value = 0.101069
$ws.Range("E44").Value = 905652369
$ws.Range("F44").Value = 59615397
$ws.Range("G44").Value = 5.33985
$ws.Range("B45").Value = "BSCX"
$ws.Range("C45").Value = "BSCEX"
$ws.Range("D45").Value = 236.13
$ws.Range("E45").Value = 900932790
$ws.Range("F45").Value = 1232324
$ws.Range("G45").Value = -0.91946
$ws.Range("B46").Value = "RPL"
$ws.Range("C46").Value = "Rocket Pool"
$ws.Range("D46").Value = 44.22
$ws.Range("E46").Value = 859118106
$ws.Range("F46").Value = 3577300
$ws.Range("G46").Value = 0.50688
$ws.Range("B47").Value = "APE"
$ws.Range("C47").Value = "ApeCoin"
$ws.Range("D47").Value = 2.33
$ws.Range("E47").Value = 855911646
$ws.Range("F47").Value = 176495193
$ws.Range("G47").Value = -0.10397
$ws.Range("B48").Value = "RETH"
$ws.Range("C48").Value = "Rocket Pool ETH"
$ws.Range("D48").Value = 1883.2
$ws.Range("E48").Value = 827072489
$ws.Range("F48").Value = 1921003
$ws.Range("G48").Value = 0.53583
$ws.Range("B49").Value = "ALGO"
$ws.Range("C49").Value = "Algorand"
$ws.Range("D49").Value = 0.108853
$ws.Range("E49").Value = 786778661
$ws.Range("F49").Value = 53314486
$ws.Range("G49").Value = 1.72047
$ws.Range("B50").Value = "EGLD"
$ws.Range("C50").Value = "MultiversX"
$ws.Range("D50").Value = 30.3
$ws.Range("E50").Value = 771665546
$ws.Range("F50").Value = 15238647
$ws.Range("G50").Value = 0.94648
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "EOS"
$ws.Range("D51").Value = 0.685325
$ws.Range("E51").Value = 757082007
$ws.Range("F51").Value = 162832180
$ws.Range("G51").Value = 0.86179
